# Add a new "ABC USA Inc." company row (row 3) to the Companies sheet,
# duplicating the existing "Origin USA Inc." row (row 2) and its three
# hyperlinks, then update A3 with the new company name and move the
# active selection to A10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Companies")

# Duplicate row 2 (values + styles + number formats) into a new row 3.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# Overwrite the company name in the new row.
$ws.Range("A3").Value = "ABC USA Inc."

# Re-create the hyperlinks that row 2 has, pointing at the same targets,
# for the corresponding cells in the new row 3.
$ws.Hyperlinks.Add($ws.Range("N3"), "https://www.facebook.com/OriginBifoldsUSA")
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.originbifolds.com/")
$ws.Hyperlinks.Add($ws.Range("O3"), "https://www.linkedin.com/company/origin-frames")

# Move / update the selection shown when the sheet is reopened.
$ws.Range("A10").Select()
